# Daily refresh of the cryptos price table (Coin/Link/Price/Volume(1h)).
# Updates the Price (D) and Volume(1h) (E) columns for each ranked coin,
# plus the two rows whose ranking swapped places (Frax <-> ImmutableX).
# Numeric-looking Price strings are apostrophe-prefixed so Excel keeps
# them as literal text (matching the source data) instead of coercing
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.122.28'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '1.852.69'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'0.6952"
$ws.Range("E5").Value = '  -4.97%  '
$ws.Range("D6").Value = "'238.61"
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = "'0.07631"
$ws.Range("E8").Value = '  +7.44%  '
$ws.Range("D9").Value = "'0.3032"
$ws.Range("E9").Value = '  -3.16%  '
$ws.Range("D10").Value = "'23.38"
$ws.Range("E10").Value = '  -4.32%  '
$ws.Range("D11").Value = "'0.08131"
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("D12").Value = '1.856.20'
$ws.Range("E12").Value = '  +1.72%  '
$ws.Range("D13").Value = "'0.7266"
$ws.Range("E13").Value = '  -2.66%  '
$ws.Range("D14").Value = "'5.235"
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("E15").Value = '  -3.68%  '
$ws.Range("D16").Value = '29.119.78'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = "'5.786"
$ws.Range("E17").Value = '  -3.95%  '
$ws.Range("D18").Value = "'0.000007749"
$ws.Range("E18").Value = '  -0.72%  '
$ws.Range("D19").Value = "'13.17"
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("D20").Value = "'236.62"
$ws.Range("E20").Value = '  -4.79%  '
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '2.096.23'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = "'7.610"
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("D25").Value = "'8.992"
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = "'161.44"
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = "'0.1447"
$ws.Range("E27").Value = '  -5.65%  '
$ws.Range("D28").Value = "'18.08"
$ws.Range("E28").Value = '  -2.51%  '
$ws.Range("D29").Value = "'1.981"
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").Value = "'1.406"
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = "'1.492"
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("D34").Value = "'0.05228"
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("E35").Value = '  -3.50%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'0.7014"
$ws.Range("E36").Value = '  -6.98%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = "'1.008"
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("E39").Value = '  -3.89%  '
$ws.Range("D40").Value = "'2.678"
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("D41").Value = "'0.9300"
$ws.Range("E41").Value = '  +7.55%  '
$ws.Range("D42").Value = "'6.028"
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = '1.083.82'
$ws.Range("E43").Value = '  +3.47%  '
$ws.Range("D44").Value = "'0.4270"
$ws.Range("E44").Value = '  -4.47%  '
$ws.Range("D45").Value = "'70.40"
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = "'103.17"
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").Value = "'1.780"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").Value = '1.992.42'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").Value = "'9.221"
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").Value = "'6.993"
$ws.Range("E51").Value = '  -6.52%  '
